$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format so numeric-looking strings (e.g. "1.000", "93.40")
# keep their exact literal representation instead of being coerced to
# numbers (which would drop trailing zeros), matching the source inlineStr cells.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "28.306.49"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.858.85"
$ws.Range("E3").Value = "  +2.81%  "
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").Value = "310.18"
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("D6").Value = "0.9978"
$ws.Range("E6").Value = "  -0.46%  "
$ws.Range("D7").Value = "0.4982"
$ws.Range("E7").Value = "  -3.28%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "0.3965"
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.1010"
$ws.Range("E9").Value = "  +28.88%  "
$ws.Range("D10").Value = "1.117"
$ws.Range("E10").Value = "  +0.48%  "
$ws.Range("D11").Value = "41.23"
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("D12").Value = "6.468"
$ws.Range("E12").Value = "  +1.95%  "
$ws.Range("D13").Value = "20.77"
$ws.Range("E13").Value = "  +1.74%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.848.72"
$ws.Range("E14").Value = "  +2.24%  "
$ws.Range("B15").Value = "BinanceUSD"
$ws.Range("C15").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D15").Value = "0.9973"
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("D16").Value = "7.387"
$ws.Range("E16").Value = "  +1.10%  "
$ws.Range("E17").Value = "  +6.07%  "
$ws.Range("D18").Value = "93.40"
$ws.Range("E18").Value = "  +0.89%  "
$ws.Range("D19").Value = "0.06655"
$ws.Range("E19").Value = "  +1.44%  "
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("D21").Value = "17.35"
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").Value = "6.050"
$ws.Range("E22").Value = "  +0.62%  "
$ws.Range("D23").Value = "28.398.97"
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("D24").Value = "11.32"
$ws.Range("E24").Value = "  +1.75%  "
$ws.Range("D25").Value = "2.248"
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "21.04"
$ws.Range("E26").Value = "  +2.55%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "2.465"
$ws.Range("E27").Value = "  +2.14%  "
$ws.Range("B28").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C28").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D28").Value = "2.044.63"
$ws.Range("E28").Value = "  +1.25%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "157.48"
$ws.Range("E29").Value = "  -2.21%  "
$ws.Range("D30").Value = "127.72"
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("D31").Value = "0.1052"
$ws.Range("E31").Value = "  -4.36%  "
$ws.Range("D32").Value = "1.053"
$ws.Range("E32").Value = "  -1.02%  "
$ws.Range("D33").Value = "5.632"
$ws.Range("E33").Value = "  +1.21%  "
$ws.Range("D34").Value = "3.596"
$ws.Range("E34").Value = "  -2.06%  "
$ws.Range("D35").Value = "0.06792"
$ws.Range("E35").Value = "  -5.46%  "
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").Value = "9.084"
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.02373"
$ws.Range("E37").Value = "  +1.15%  "
$ws.Range("D38").Value = "0.2159"
$ws.Range("E38").Value = "  -0.94%  "
$ws.Range("D39").Value = "5.024"
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("D40").Value = "11.48"
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("D41").Value = "0.6268"
$ws.Range("E41").Value = "  +1.26%  "
$ws.Range("D42").Value = "1.176"
$ws.Range("E42").Value = "  +1.89%  "
$ws.Range("D43").Value = "0.9978"
$ws.Range("E43").Value = "  -0.38%  "
$ws.Range("D44").Value = "13.30"
$ws.Range("E44").Value = "  +0.46%  "
$ws.Range("D45").Value = "0.5976"
$ws.Range("E45").Value = "  -0.16%  "
$ws.Range("D46").Value = "3.698"
$ws.Range("E46").Value = "  -1.35%  "
$ws.Range("D47").Value = "1.279"
$ws.Range("E47").Value = "  -2.01%  "
$ws.Range("D48").Value = "124.61"
$ws.Range("E48").Value = "  -0.50%  "
$ws.Range("D49").Value = "1.964"
$ws.Range("E49").Value = "  +2.22%  "
$ws.Range("D50").Value = "1.189"
$ws.Range("E50").Value = "  -2.48%  "
$ws.Range("E51").Value = "  +4.55%  "
